$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B4").NumberFormat = "@"

$ws.Range("A2").Value = "关于准予撤销上海市嘉定区朱家桥邮政支局的公告"
$ws.Range("B2").Value = "2025-12-05"
$ws.Range("C2").Value = "https://sh.spb.gov.cn/shsyzglj/c100057/c100058/202512/b487344d328e4d2fa163d4fe9a0fe502.shtml"

$ws.Range("A3").Value = "关于准予撤销上海市徐汇区柳州路邮政所的公告"
$ws.Range("B3").Value = "2025-12-02"
$ws.Range("C3").Value = "https://sh.spb.gov.cn/shsyzglj/c100057/c100058/202512/be58981880de42c7822366e7faabd2cb.shtml"

$ws.Range("A4").Value = "关于准予撤销上海市松江区达丰邮政所的公告"
$ws.Range("B4").Value = "2025-10-31"
$ws.Range("C4").Value = "https://sh.spb.gov.cn/shsyzglj/c100057/c100058/202510/e5a44a5099d1476fa0e479b321267ac3.shtml"
